# DW+P5 - Modele plan tests acceptation
# Adds a new "Gestion des erreurs" column (E) between the existing
# "Résultat attendu" (D) and "Résultat observé" (E->F) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Make room for the new column: duplicate column E (values + styles)
#        into the new column F. This is what shifts the former "Résultat
#        observé" content one column to the right while keeping its
#        original formatting intact.
$ws.Range("E1:E22").Copy($ws.Range("F1:F22"))

# --- 2) Column F should have the same width as column E.
$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth

# --- 3) Fill in the new "Gestion des erreurs" column (E) content.
$ws.Range("E1").Value = "Gestion des erreurs"
$ws.Range("E2").Value = "Si erreur, afficher un message. Proposer d'envoyer un message au responsable du site."
$ws.Range("E4").Value = "Si aucune couleur specifié, afficher un message d'erreur pour la correction des champs."
$ws.Range("E5").Value = "Si la quantité specifié est hors limites, afficher un message d'erreur pour la correction des champs.."
$ws.Range("E6").Value = "Si erreur, afficher un message. Proposer d'envoyer un message au responsable du site."
$ws.Range("E7").Value = "Si erreur, afficher un message. Proposer d'envoyer un message au responsable du site."
$ws.Range("E8").Value = "NA"
$ws.Range("E9").Value = "Si erreur, afficher un message. Proposer d'envoyer un message au responsable du site."
$ws.Range("E10").ClearContents()
$ws.Range("E11").Value = "Si erreur, afficher un message. Proposer d'envoyer un message au responsable du site."
$ws.Range("E12").Value = "Si un champ n'est pas valide, afficher un message d'erreur pour la correction des champs.."
$ws.Range("E13").Value = "Si erreur, afficher un message. Proposer d'envoyer un message au responsable du site."
$ws.Range("E14").Value = "Si erreur, afficher un message. Proposer d'envoyer un message au responsable du site."
# E3 keeps its original "OK" value/style (same as the new F3) - no change needed.

# --- 4) Rows whose new "Gestion des erreurs" text wraps onto more lines
#        grow to match the taller rows used elsewhere in the sheet.
$ws.Rows(6).RowHeight = 65.25
$ws.Rows(7).RowHeight = 65.25
$ws.Rows(11).RowHeight = 65.25
$ws.Rows(13).RowHeight = 65.25
$ws.Rows(14).RowHeight = 65.25

# --- 5) View state: scrolled down a little further, selection parked on
#        the last "Gestion des erreurs" cell that was just edited.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("E14").Select()
